# Updates cryptos list values (Price and Volume(1h) columns) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.320.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6204"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07338"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2875"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6595"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.226"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.309.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "236.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.184"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.391"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1329"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06878"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.474"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.002"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.918"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.150"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.741"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6790"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.582"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01814"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.230.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.634"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9442"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.988.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.680"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.861"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.784"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1123"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.62%  "
